$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("(BD)", $true, $false, $false, $false, $false, $true, 1, $false,
              "(BD) + 2/3 din arhitectura aplicației", 2)
